$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the rate summary text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cell = $ws1.Range("A1")
$text = $cell.Value2

$old1 = "1000 Bs = 6.98 = 28051.25 pesos"
$new1 = "1000 Bs = 7.12 = 28511.6 pesos"
$old2 = "28051.25 pesos = 6.93 = 963.23 Bs"
$new2 = "28511.6 pesos = 7.09 = 951.84 Bs"

$text = $text.Replace($old1, $new1)
$text = $text.Replace($old2, $new2)
$cell.Value = $text

# --- Sheet "tasas": update the rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 140.399
$ws2.Range("O10").Value = 4003
$ws2.Range("N12").Value = 4019.9
$ws2.Range("O12").Value = 134.201
